$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "Leadership and Line Managers" with "Leadership" everywhere it is used
# (column C, Category). Using Find/Replace across the used range mirrors what a
# user would do in Excel (Ctrl+H) to rename a category label.
$usedRange = $ws.UsedRange
$usedRange.Replace("Leadership and Line Managers", "Leadership", 1, 1, $false, $false, $false, $false)

# Update the active selection left by the user to C4
$ws.Range("C4").Select()
